$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Columns("C").Insert()
